# Natmi following Dr Hou advice
# Update LR-pairs_lrc2p/Bmp4-Bmpr2 sheet: rebuild sending/target cluster combos
# for FAPs / sCs / ECs (3x3 = 9 rows) with refreshed expression stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp4"
$ws.Cells.Item(2, 3).Value = "Bmpr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 5.448315
$ws.Cells.Item(2, 8).Value = 16.344945
$ws.Cells.Item(2, 9).Value = 0.1199618029178375
$ws.Cells.Item(2, 10).Value = 0.1199618029178375
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 40.70766766666667
$ws.Cells.Item(2, 14).Value = 122.123003
$ws.Cells.Item(2, 15).Value = 0.3776398983502007
$ws.Cells.Item(2, 16).Value = 0.3776398983502007
$ws.Cells.Item(2, 17).Value = 221.788196363315
$ws.Cells.Item(2, 18).Value = 1996.093767269835
$ws.Cells.Item(2, 19).Value = 0.04530236305979898
$ws.Cells.Item(2, 20).Value = 0.04530236305979898

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp4"
$ws.Cells.Item(3, 3).Value = "Bmpr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 5.448315
$ws.Cells.Item(3, 8).Value = 16.344945
$ws.Cells.Item(3, 9).Value = 0.1199618029178375
$ws.Cells.Item(3, 10).Value = 0.1199618029178375
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 39.715023
$ws.Cells.Item(3, 14).Value = 119.145069
$ws.Cells.Item(3, 15).Value = 0.3684312589831062
$ws.Cells.Item(3, 16).Value = 0.3684312589831062
$ws.Cells.Item(3, 17).Value = 216.379955536245
$ws.Cells.Item(3, 18).Value = 1947.419599826205
$ws.Cells.Item(3, 19).Value = 0.04419767807890215
$ws.Cells.Item(3, 20).Value = 0.04419767807890215

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bmp4"
$ws.Cells.Item(4, 3).Value = "Bmpr2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 5.448315
$ws.Cells.Item(4, 8).Value = 16.344945
$ws.Cells.Item(4, 9).Value = 0.1199618029178375
$ws.Cells.Item(4, 10).Value = 0.1199618029178375
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 27.37224266666666
$ws.Cells.Item(4, 14).Value = 82.116728
$ws.Cells.Item(4, 15).Value = 0.253928842666693
$ws.Cells.Item(4, 16).Value = 0.253928842666693
$ws.Cells.Item(4, 17).Value = 149.13260030444
$ws.Cells.Item(4, 18).Value = 1342.19340273996
$ws.Cells.Item(4, 19).Value = 0.03046176177913641
$ws.Cells.Item(4, 20).Value = 0.03046176177913641

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bmp4"
$ws.Cells.Item(5, 3).Value = "Bmpr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 31.28251333333334
$ws.Cells.Item(5, 8).Value = 93.84754000000001
$ws.Cells.Item(5, 9).Value = 0.6887829905701046
$ws.Cells.Item(5, 10).Value = 0.6887829905701045
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 40.70766766666667
$ws.Cells.Item(5, 14).Value = 122.123003
$ws.Cells.Item(5, 15).Value = 0.3776398983502007
$ws.Cells.Item(5, 16).Value = 0.3776398983502007
$ws.Cells.Item(5, 17).Value = 1273.438156551402
$ws.Cells.Item(5, 18).Value = 11460.94340896262
$ws.Cells.Item(5, 19).Value = 0.2601119385442416
$ws.Cells.Item(5, 20).Value = 0.2601119385442415

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bmp4"
$ws.Cells.Item(6, 3).Value = "Bmpr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 31.28251333333334
$ws.Cells.Item(6, 8).Value = 93.84754000000001
$ws.Cells.Item(6, 9).Value = 0.6887829905701046
$ws.Cells.Item(6, 10).Value = 0.6887829905701045
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 39.715023
$ws.Cells.Item(6, 14).Value = 119.145069
$ws.Cells.Item(6, 15).Value = 0.3684312589831062
$ws.Cells.Item(6, 16).Value = 0.3684312589831062
$ws.Cells.Item(6, 17).Value = 1242.38573653114
$ws.Cells.Item(6, 18).Value = 11181.47162878026
$ws.Cells.Item(6, 19).Value = 0.2537691843818926
$ws.Cells.Item(6, 20).Value = 0.2537691843818926

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Bmp4"
$ws.Cells.Item(7, 3).Value = "Bmpr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 31.28251333333334
$ws.Cells.Item(7, 8).Value = 93.84754000000001
$ws.Cells.Item(7, 9).Value = 0.6887829905701046
$ws.Cells.Item(7, 10).Value = 0.6887829905701045
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 27.37224266666666
$ws.Cells.Item(7, 14).Value = 82.116728
$ws.Cells.Item(7, 15).Value = 0.253928842666693
$ws.Cells.Item(7, 16).Value = 0.253928842666693
$ws.Cells.Item(7, 17).Value = 856.2725461832356
$ws.Cells.Item(7, 18).Value = 7706.45291564912
$ws.Cells.Item(7, 19).Value = 0.1749018676439704
$ws.Cells.Item(7, 20).Value = 0.1749018676439704

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Bmp4"
$ws.Cells.Item(8, 3).Value = "Bmpr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.686253333333333
$ws.Cells.Item(8, 8).Value = 26.05876
$ws.Cells.Item(8, 9).Value = 0.1912552065120579
$ws.Cells.Item(8, 10).Value = 0.1912552065120579
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 40.70766766666667
$ws.Cells.Item(8, 14).Value = 122.123003
$ws.Cells.Item(8, 15).Value = 0.3776398983502007
$ws.Cells.Item(8, 16).Value = 0.3776398983502007
$ws.Cells.Item(8, 17).Value = 353.5971139618089
$ws.Cells.Item(8, 18).Value = 3182.37402565628
$ws.Cells.Item(8, 19).Value = 0.07222559674616019
$ws.Cells.Item(8, 20).Value = 0.07222559674616019

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Bmp4"
$ws.Cells.Item(9, 3).Value = "Bmpr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 8.686253333333333
$ws.Cells.Item(9, 8).Value = 26.05876
$ws.Cells.Item(9, 9).Value = 0.1912552065120579
$ws.Cells.Item(9, 10).Value = 0.1912552065120579
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 39.715023
$ws.Cells.Item(9, 14).Value = 119.145069
$ws.Cells.Item(9, 15).Value = 0.3684312589831062
$ws.Cells.Item(9, 16).Value = 0.3684312589831062
$ws.Cells.Item(9, 17).Value = 344.97475091716
$ws.Cells.Item(9, 18).Value = 3104.77275825444
$ws.Cells.Item(9, 19).Value = 0.07046439652231147
$ws.Cells.Item(9, 20).Value = 0.07046439652231147

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Bmp4"
$ws.Cells.Item(10, 3).Value = "Bmpr2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 8.686253333333333
$ws.Cells.Item(10, 8).Value = 26.05876
$ws.Cells.Item(10, 9).Value = 0.1912552065120579
$ws.Cells.Item(10, 10).Value = 0.1912552065120579
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 27.37224266666666
$ws.Cells.Item(10, 14).Value = 82.116728
$ws.Cells.Item(10, 15).Value = 0.253928842666693
$ws.Cells.Item(10, 16).Value = 0.253928842666693
$ws.Cells.Item(10, 17).Value = 237.7622341041422
$ws.Cells.Item(10, 18).Value = 2139.86010693728
$ws.Cells.Item(10, 19).Value = 0.04856521324358624
$ws.Cells.Item(10, 20).Value = 0.04856521324358624
